$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vokabeltrainer Sprint 1")
$ws.Activate()
$ws.Range("D5").Value = 0
$ws.Range("D6").Select()
